$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace with new Pertamina/Kilang news item ---
$ws.Range("A2").Value = "Purbaya Cap Pertamina Malas-malasan Bangun KilangMONEY01/10/2025"

# Keep the date column as literal text (matches original inlineStr text),
# not an auto-converted Excel date serial number.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2025-10-01"
$ws.Range("B2").Style = "Normal"

$ws.Range("D2").Value = "https://money.kompas.com/read/2025/10/01/080353726/purbaya-cap-pertamina-malas-malasan-bangun-kilang"
$ws.Range("E2").Value = "purbaya"

# --- Row 3: replace with new Pertamina kilang criticism news item ---
$ws.Range("A3").Value = "Kritik Menkeu Purbaya, Pertamina Gagal Realisasikan Kilang Baru Setelah 7 TahunPROV01/10/2025"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2025-10-01"
$ws.Range("B3").Style = "Normal"

$ws.Range("D3").Value = "https://www.kompas.com/sumatera-selatan/read/2025/10/01/051500288/kritik-menkeu-purbaya-pertamina-gagal-realisasikan-kilang-baru"
$ws.Range("E3").Value = "purbaya"

# --- Remove rows 4 and 5 entirely (dimension shrinks to A1:E3) ---
$ws.Range("A4:E5").Delete()
